# Add a caption textbox below the "Overall Performance Metrics" picture on
# slide 18, describing the benchmark hardware / units used in the table.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)

# Shape position/size properties in this COM layer are expressed in points
# (standard PowerPoint VBA semantics: 1 pt = 12700 EMU), so convert the
# target EMU values from the OOXML diff back to points here. Passing the
# converted points straight into AddTextbox(...) round-trips back to the
# exact EMU integers on save.
$emuPerPt = 12700

$left   = 4007627 / $emuPerPt
$top    = 2642911 / $emuPerPt
$width  = 4176744 / $emuPerPt
$height = 646331  / $emuPerPt

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "TextBox 2"

# Shape formatting: no fill, word-wrap with "resize shape to fit text".
$tb.Fill.Visible = 0
$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 1

# Build the run-split text exactly like the source document: the German
# words are typed as separate runs (this is also how PowerPoint's spell
# checker ends up flagging them individually).
$tr = $tb.TextFrame.TextRange
$tr.Text = "AMD R5 3600, DDR4-3200, GTX 1070, "
[void]$tr.InsertAfter("Zeiten")
[void]$tr.InsertAfter(" in ")
[void]$tr.InsertAfter("ms")
